$wb = $excel.ActiveWorkbook

# 1. Rename "Hoja1" -> "General"
$ws1 = $wb.Worksheets.Item("Hoja1")
$ws1.Name = "General"

$ws2 = $wb.Worksheets.Item("Agregar jugadores")

# 2. Mark existing rows as "ok" in the Status column (D)
$ws1.Range("D3").Value = "ok"
$ws1.Range("D4").Value = "ok"
$ws1.Range("D6").Value = "ok"
$ws1.Range("D9").Value = "ok"
$ws1.Range("D10").Value = "ok"
$ws1.Range("D14").Value = "ok"
$ws1.Range("D15").Value = "ok"
$ws1.Range("D16").Value = "ok"
$ws1.Range("D17").Value = "ok"

# 3. Update "Agregar jugadores" sheet
$ws2.Range("C3").Value = "*"
$ws2.Range("D4").Value = "¿Configurar qué?"

# 4. Add new backlog item as row 18
$ws1.Range("A18").Value = 17
$ws1.Range("B18").Value = "Agregar más preguntas"
$ws1.Range("C18").Value = "Y sus respuestas. También determinar el orden de los niveles."

# 5. Add description for row 11 (Cliente de prueba)
$ws1.Range("C11").Value = "Que muestre lo que se le preguntó y conteste algún valor hardcodeado"

# 6. Update selections - select sheet2's cell first, then sheet1's so the
#    "General" sheet remains the active/selected tab.
$ws2.Range("B6").Select()
$ws1.Range("D18").Select()
